$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) - reorder/rename the distractor columns so that
# n distractor = n targets (bedrooms columns moved next to their pair)
$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("B1").Value = "bedrooms_1"
$ws.Range("C1").Value = "kitchens_1"
$ws.Range("D1").Value = "living_rooms_2"
$ws.Range("E1").Value = "bedrooms_2"
$ws.Range("F1").Value = "kitchens_2"

# Update the data rows (2-7) with the new block-order one-hot values
$data = @(
    @(0,0,0,0,1,0),
    @(0,0,0,0,0,1),
    @(1,0,0,0,0,0),
    @(0,1,0,0,0,0),
    @(0,0,1,0,0,0),
    @(0,0,0,1,0,0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
    }
}
